$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -5
